# Add three new work-log entries (rows 28-30) documenting the car prefab
# update and the new Grab blend tree work, mirroring the formatting already
# used by the existing rows in this time-log table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the time-column formatting (style used by B27/C27) for the new B/C cells.
$ws.Range("B28:C29").NumberFormat = $ws.Range("B27").NumberFormat
$ws.Range("B28:C29").HorizontalAlignment = $ws.Range("B27").HorizontalAlignment
$ws.Range("B30").NumberFormat = $ws.Range("B27").NumberFormat
$ws.Range("B30").HorizontalAlignment = $ws.Range("B27").HorizontalAlignment

# Match the description-column formatting (style used by D27) for the new D cells.
$ws.Range("D28:D30").HorizontalAlignment = $ws.Range("D27").HorizontalAlignment
$ws.Range("D28:D30").VerticalAlignment = $ws.Range("D27").VerticalAlignment
$ws.Range("D28:D30").WrapText = $ws.Range("D27").WrapText

# Row 28: 20:15 - 20:45, updating the car prefab.
$ws.Range("B28").Value2 = 0.84375
$ws.Range("C28").Value2 = 0.86458333333333337
$ws.Range("D28").Value = "Updating the prefab to the new car."

# Row 29: 21:00 - 21:10, researching animation layers for the hands.
$ws.Range("B29").Value2 = 0.875
$ws.Range("C29").Value2 = 0.88194444444444453
$ws.Range("D29").Value = "Reseraching Animation layers for the hands."

# Row 30: starts at 21:10, still working - importing a timer from another project.
$ws.Range("B30").Value2 = 0.88194444444444453
$ws.Range("D30").Value = "Importing timer from my other project."

$ws.Range("D30").Select() | Out-Null
